$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.869.46"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "2.105.51"
$ws.Range("E3").Value = "  +3.36%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'228.71"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("D7").Value = "'60.44"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.385"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "2.416.33"
$ws.Range("E12").Value = "  +3.33%  "
$ws.Range("E13").Value = "  +3.81%  "
$ws.Range("D14").Value = "'22.21"
$ws.Range("E14").Value = "  +5.72%  "
$ws.Range("D15").Value = "'0.797"
$ws.Range("E15").Value = "  +3.42%  "
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "2.112.59"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "38.799.03"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("D19").Value = "'72.09"
$ws.Range("E19").Value = "  +3.89%  "
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").Value = "'226.15"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").Value = "'170.78"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "'9.53"
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("E28").Value = "  +6.44%  "
$ws.Range("D29").Value = "'1.40"
$ws.Range("E29").Value = "  +9.45%  "
$ws.Range("D30").Value = "'19.22"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  +4.62%  "
$ws.Range("D33").Value = "'4.75"
$ws.Range("E33").Value = "  +6.25%  "
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").Value = "'2.40"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'101.69"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.545.72"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").Value = "'0.0930"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").Value = "'7.71"
$ws.Range("E46").Value = "  +9.07%  "
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").Value = "2.302.04"
$ws.Range("E51").Value = "  +3.28%  "